{"js": "// Update the date line and the 25 division-equation cells in the table\n// to the new values, preserving existing run/paragraph formatting by\n// replacing text in-place (via Range.insertText with InsertLocation.Replace)\n// rather than rebuilding paragraphs/runs.\n\nconst body = context.document.body;\n\n// --- 1) Update the date heading paragraph ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.getRange().insertText(\"2024-12-31 Tuesday\", Word.InsertLocation.replace);\n\n// --- 2) Update the table of division problems ---\nconst table = body.tables.getFirst();\n\n// Row indices that contain data (every 4th row: data row followed by 3 blank rows)\nconst dataRows = [0, 4, 8, 12, 16];\n\n// New values per data row, in column order (0-4)\nconst newValues = [\n  [\"79\u00f75=15, 4\", \"50\u00f74=12, 2\", \"54\u00f73=18, 0\", \"40\u00f76=6, 4\", \"54\u00f77=7, 5\"],\n  [\"47\u00f77=6, 5\", \"32\u00f73=10, 2\", \"38\u00f74=9, 2\", \"27\u00f76=4, 3\", \"75\u00f72=37, 1\"],\n  [\"59\u00f78=7, 3\", \"26\u00f77=3, 5\", \"47\u00f75=9, 2\", \"66\u00f75=13, 1\", \"74\u00f79=8, 2\"],\n  [\"91\u00f79=10, 1\", \"46\u00f78=5, 6\", \"89\u00f77=12, 5\", \"28\u00f73=9, 1\", \"44\u00f75=8, 4\"],\n  [\"41\u00f78=5, 1\", \"38\u00f76=6, 2\", \"58\u00f76=9, 4\", \"51\u00f77=7, 2\", \"78\u00f76=13, 0\"],\n];\n\nfor (let i = 0; i < dataRows.length; i++) {\n  const rowIndex = dataRows[i];\n  for (let col = 0; col < 5; col++) {\n    const cell = table.getCell(rowIndex, col);\n    const cellParagraph = cell.body.paragraphs.getFirst();\n    cellParagraph.getRange().insertText(newValues[i][col], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 division-equation cells in the table\n# to the new values, preserving existing run/paragraph formatting by\n# setting Range.Text in place (keeps rPr/pPr untouched).\n\n$d = $word.ActiveDocument\n\n# --- 1) Update the date heading paragraph (first paragraph in the document) ---\n$d.Paragraphs.Item(1).Range.Text = \"2024-12-31 Tuesday\"\n\n# --- 2) Update the table of division problems ---\n$t = $d.Tables.Item(1)\n\n# Row indices (1-based) that contain data (every 4th row: data row followed by 3 blank rows)\n$dataRows = @(1, 5, 9, 13, 17)\n\n# New values per data row, in column order (1-5)\n$newValues = @(\n  @(\"79\u00f75=15, 4\", \"50\u00f74=12, 2\", \"54\u00f73=18, 0\", \"40\u00f76=6, 4\", \"54\u00f77=7, 5\"),\n  @(\"47\u00f77=6, 5\", \"32\u00f73=10, 2\", \"38\u00f74=9, 2\", \"27\u00f76=4, 3\", \"75\u00f72=37, 1\"),\n  @(\"59\u00f78=7, 3\", \"26\u00f77=3, 5\", \"47\u00f75=9, 2\", \"66\u00f75=13, 1\", \"74\u00f79=8, 2\"),\n  @(\"91\u00f79=10, 1\", \"46\u00f78=5, 6\", \"89\u00f77=12, 5\", \"28\u00f73=9, 1\", \"44\u00f75=8, 4\"),\n  @(\"41\u00f78=5, 1\", \"38\u00f76=6, 2\", \"58\u00f76=9, 4\", \"51\u00f77=7, 2\", \"78\u00f76=13, 0\")\n)\n\nfor ($i = 0; $i -lt $dataRows.Length; $i++) {\n  $row = $dataRows[$i]\n  for ($col = 1; $col -le 5; $col++) {\n    $cell = $t.Cell($row, $col)\n    $cell.Range.Text = $newValues[$i][$col - 1]\n  }\n}\n"}
